$p = $ppt.ActivePresentation

# --- Slide 17 ("The centering Algorithm"): fix "Cb(Un) = forward center..."
#     which should read "Cp(Un) = forward center..." (typo fix: b -> p)
$s17 = $p.Slides.Item(17)
$shape17 = $s17.Shapes.Item(2)
$tr17 = $shape17.TextFrame.TextRange
$found17 = $tr17.Find("Cb(Un) = forward")
if ($found17 -ne $null) {
    $bPos = $found17.Start + 1
    $tr17.Characters($bPos, 1).Text = "p"
}

# --- Slide 22 ("Further reading"): merge the split " " / "21.1-21.3, " /
#     "21.5-21.6" runs back into a single run reading " 21.1-21.3, 21.5-21.6"
$s22 = $p.Slides.Item(22)
$shape22 = $s22.Shapes.Item(2)
$tr22 = $shape22.TextFrame.TextRange
$found22 = $tr22.Find("21.1-21.3")
if ($found22 -ne $null) {
    $mergeStart = $found22.Start - 1
    $tr22.Characters($mergeStart, 21).Text = " 21.1-21.3, 21.5-21.6"
}
